$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = 73962
$ws.Range("B2").Value = "Mariana Pacheco"
$ws.Range("C2").Value = "TI"
$ws.Range("F2").Value = 45100
$ws.Range("G2").Value = 5014.94

# Row 3
$ws.Range("A3").Value = 6229
$ws.Range("B3").Value = "Cauã Vasconcelos"
$ws.Range("C3").Value = "Engenharia"
$ws.Range("D3").Value = "Consulta medica"
$ws.Range("F3").Value = 45096
$ws.Range("G3").Value = 3903.68

# Row 4
$ws.Range("A4").Value = 86469
$ws.Range("B4").Value = "Isabella Pereira"
$ws.Range("D4").Value = "Outros"
$ws.Range("E4").Value = 6
$ws.Range("F4").Value = 45101
$ws.Range("G4").Value = 4295.07

# Row 5
$ws.Range("A5").Value = 12071
$ws.Range("B5").Value = "Dr. Levi Viana"
$ws.Range("D5").Value = "Doenca"
$ws.Range("E5").Value = 7
$ws.Range("F5").Value = 45095
$ws.Range("G5").Value = 8412.51

# Row 6
$ws.Range("A6").Value = 80174
$ws.Range("B6").Value = "Rafael Gonçalves"
$ws.Range("C6").Value = "Vendas"
$ws.Range("D6").Value = "Doenca"
$ws.Range("E6").Value = 8
$ws.Range("F6").Value = 45091
$ws.Range("G6").Value = 6560.05

# Row 7
$ws.Range("A7").Value = 60741
$ws.Range("B7").Value = "Hellena Vasconcelos"
$ws.Range("C7").Value = "Engenharia"
$ws.Range("D7").Value = "Consulta medica"
$ws.Range("E7").Value = 6
$ws.Range("F7").Value = 45095
$ws.Range("G7").Value = 3293.6

# Row 8
$ws.Range("A8").Value = 58680
$ws.Range("B8").Value = "Alice Vasconcelos"
$ws.Range("C8").Value = "Financeiro"
$ws.Range("D8").Value = "Doenca"
$ws.Range("E8").Value = 2
$ws.Range("F8").Value = 45086
$ws.Range("G8").Value = 3085.02

# Row 9
$ws.Range("A9").Value = 4483
$ws.Range("B9").Value = "Mateus da Rocha"
$ws.Range("C9").Value = "Juridico"
$ws.Range("D9").Value = "Doenca"
$ws.Range("E9").Value = 6
$ws.Range("F9").Value = 45098
$ws.Range("G9").Value = 7799.98

# Row 10
$ws.Range("A10").Value = 39374
$ws.Range("B10").Value = "Eloah Freitas"
$ws.Range("D10").Value = "Viagem de negocios"
$ws.Range("F10").Value = 45088
$ws.Range("G10").Value = 5001.39

# Row 11
$ws.Range("A11").Value = 20484
$ws.Range("B11").Value = "Ana Júlia Novais"
$ws.Range("C11").Value = "Financeiro"
$ws.Range("E11").Value = 5
$ws.Range("F11").Value = 45101
$ws.Range("G11").Value = 4380.24
